# Add "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell H1, styled the same as the other header cells (B1:G1 use style index 1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Save values for each data row (H2:H23), taken from the source data.
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
